# Update the SEM path-coefficient labels on the single figure slide to the
# refreshed model estimates (fresh SEM table layout).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$updates = @{
    19 = "-0.079"
    20 = "-0.022"
    24 = "0.026"
    25 = "-0.059"
    26 = "-0.039"
    27 = "0.020"
    28 = "-0.173"
    29 = "-0.081"
    32 = "-0.902"
    34 = "-0.774"
    36 = "-0.602"
    38 = "-0.540"
    42 = "0.527"
}

foreach ($idx in $updates.Keys) {
    $shape = $s.Shapes.Item([int]$idx)
    $shape.TextFrame.TextRange.Text = $updates[$idx]
}
